# Update the "取得日時" (fetch timestamp) column (A) for rows 2-14
# on the "ランサーズ" sheet from "2025-12-12 18:29:17" to
# "2025-12-12 18:37:18", as described by the commit:
#   "Append: 2025-12-12 18:37 JST"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-12 18:37:18"

for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
